# Auto-generated script: applies a bulk data refresh to the Leve profit tracking sheets.
# For each sheet, cells in columns H-N (computed market/profit figures) are updated
# to reflect newly refreshed market data. A few cells are cleared (no longer applicable)
# and a couple of previously-empty profit cells are populated.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10500
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 10500
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H38").Value = 1199.5714
$ws.Range("I38").Value = 701.3333
$ws.Range("J38").Value = 1573.25
$ws.Range("K38").Value = 2103.9999
$ws.Range("L38").Value = 4719.75
$ws.Range("M38").Value = -1731.9999
$ws.Range("N38").Value = -5463.75
$ws.Range("H43").Value = 1958.1666
$ws.Range("I43").Value = 2666.3333
$ws.Range("K43").Value = 2666.3333
$ws.Range("M43").Value = -2597.3333
$ws.Range("H48").Value = 5151.8335
$ws.Range("I48").Value = 3456
$ws.Range("K48").Value = 10368
$ws.Range("M48").Value = -10076
$ws.Range("H56").Value = 5151.8335
$ws.Range("I56").Value = 3456
$ws.Range("K56").Value = 10368
$ws.Range("M56").Value = -9834
$ws.Range("H64").Value = 9571.429
$ws.Range("I64").Value = 8612.75
$ws.Range("J64").Value = 9954.9
$ws.Range("K64").Value = 8612.75
$ws.Range("L64").Value = 9954.9
$ws.Range("M64").Value = -8364.75
$ws.Range("N64").Value = -10450.9
$ws.Range("H67").Value = 9571.429
$ws.Range("I67").Value = 8612.75
$ws.Range("J67").Value = 9954.9
$ws.Range("K67").Value = 8612.75
$ws.Range("L67").Value = 9954.9
$ws.Range("M67").Value = -7754.75
$ws.Range("N67").Value = -11670.9
$ws.Range("H76").Value = 4700.5
$ws.Range("I76").Value = 3767.3333
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 3767.3333
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -3452.3333
$ws.Range("N76").Value = -8130
$ws.Range("H79").Value = 4700.5
$ws.Range("I79").Value = 3767.3333
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 3767.3333
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -2675.3333
$ws.Range("N79").Value = -9684
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3983.4546
$ws.Range("I63").Value = 2479.8
$ws.Range("J63").Value = 5236.5
$ws.Range("K63").Value = 2479.8
$ws.Range("L63").Value = 5236.5
$ws.Range("M63").Value = -1793.8
$ws.Range("N63").Value = -6608.5
$ws.Range("H66").Value = 3983.4546
$ws.Range("I66").Value = 2479.8
$ws.Range("J66").Value = 5236.5
$ws.Range("K66").Value = 12399
$ws.Range("L66").Value = 26182.5
$ws.Range("M66").Value = -8967
$ws.Range("N66").Value = -33046.5
$ws.Range("H74").Value = 6977.5884
$ws.Range("I74").Value = 6583.2144
$ws.Range("K74").Value = 6583.2144
$ws.Range("M74").Value = -5709.2144
$ws.Range("H77").Value = 6977.5884
$ws.Range("I77").Value = 6583.2144
$ws.Range("K77").Value = 32916.072
$ws.Range("M77").Value = -28548.072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3363.5881
$ws.Range("I20").Value = 3272
$ws.Range("J20").Value = 3583.4
$ws.Range("K20").Value = 3272
$ws.Range("L20").Value = 3583.4
$ws.Range("M20").Value = -3025
$ws.Range("N20").Value = -4077.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4189.5
$ws.Range("I99").Value = 3973.5557
$ws.Range("K99").Value = 3973.5557
$ws.Range("M99").Value = -2475.5557
$ws.Range("H126").Value = 4189.5
$ws.Range("I126").Value = 3973.5557
$ws.Range("K126").Value = 11920.6671
$ws.Range("M126").Value = -9450.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1800.2858
$ws.Range("I75").Value = 1747.5
$ws.Range("J75").Value = 1821.4
$ws.Range("K75").Value = 5242.5
$ws.Range("L75").Value = 5464.200000000001
$ws.Range("M75").Value = -4244.5
$ws.Range("N75").Value = -7460.200000000001
$ws.Range("H78").Value = 1800.2858
$ws.Range("I78").Value = 1747.5
$ws.Range("J78").Value = 1821.4
$ws.Range("K78").Value = 15727.5
$ws.Range("L78").Value = 16392.6
$ws.Range("M78").Value = -10735.5
$ws.Range("N78").Value = -26376.6
$ws.Range("H97").Value = 3994
$ws.Range("I97").Value = 1999.5
$ws.Range("K97").Value = 5998.5
$ws.Range("M97").Value = -5502.5
$ws.Range("H103").Value = 2999.4
$ws.Range("J103").Value = 3249.25
$ws.Range("L103").Value = 9747.75
$ws.Range("N103").Value = -11505.75
$ws.Range("H114").Value = 3305.2
$ws.Range("I114").Value = 1831.6666
$ws.Range("J114").Value = 5515.5
$ws.Range("K114").Value = 5494.9998
$ws.Range("L114").Value = 16546.5
$ws.Range("M114").Value = -2240.9998
$ws.Range("N114").Value = -23054.5
$ws.Range("H131").Value = 3310.1538
$ws.Range("J131").Value = 4113.857
$ws.Range("L131").Value = 12341.571
$ws.Range("N131").Value = -22421.571
$ws.Range("H132").Value = 1959.75
$ws.Range("I132").Value = 1199.5
$ws.Range("J132").Value = 2213.1667
$ws.Range("K132").Value = 10795.5
$ws.Range("L132").Value = 19918.5003
$ws.Range("M132").Value = -8265.5
$ws.Range("N132").Value = -24978.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62586.375
$ws.Range("I2").Value = 100044.2
$ws.Range("K2").Value = 100044.2
$ws.Range("M2").Value = -99931.2
$ws.Range("H15").Value = 21500
$ws.Range("J15").Value = 21500
$ws.Range("L15").Value = 21500
$ws.Range("N15").Value = -22076
$ws.Range("H80").Value = 2011.6666
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 2317.5
$ws.Range("K80").Value = 1400
$ws.Range("L80").Value = 2317.5
$ws.Range("M80").Value = -402
$ws.Range("N80").Value = -4313.5
$ws.Range("H81").Value = 21500
$ws.Range("J81").Value = 21500
$ws.Range("L81").Value = 21500
$ws.Range("N81").Value = -23496
$ws.Range("H83").Value = 2011.6666
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 2317.5
$ws.Range("K83").Value = 7000
$ws.Range("L83").Value = 11587.5
$ws.Range("M83").Value = -2008
$ws.Range("N83").Value = -21571.5
$ws.Range("H84").Value = 21500
$ws.Range("J84").Value = 21500
$ws.Range("L84").Value = 64500
$ws.Range("N84").Value = -74484
$ws.Range("H102").Value = 1730
$ws.Range("I102").Value = 1730
$ws.Range("K102").Value = 1730
$ws.Range("M102").Value = -108
$ws.Range("H132").Value = 36109.938
$ws.Range("I132").Value = 44063
$ws.Range("J132").Value = 7706.143
$ws.Range("K132").Value = 132189
$ws.Range("L132").Value = 23118.429
$ws.Range("M132").Value = -129659
$ws.Range("N132").Value = -28178.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5497.727
$ws.Range("I61").Value = 3420.8333
$ws.Range("K61").Value = 3420.8333
$ws.Range("M61").Value = -3218.8333
$ws.Range("H82").Value = 2834.9412
$ws.Range("I82").Value = 1549.5
$ws.Range("K82").Value = 1549.5
$ws.Range("M82").Value = -1188.5
$ws.Range("H85").Value = 2834.9412
$ws.Range("I85").Value = 1549.5
$ws.Range("K85").Value = 1549.5
$ws.Range("M85").Value = -301.5
$ws.Range("H113").Value = 5497.727
$ws.Range("I113").Value = 3420.8333
$ws.Range("K113").Value = 3420.8333
$ws.Range("M113").Value = -1250.8333
$ws.Range("H132").Value = 2963.6365
$ws.Range("I132").Value = 2994.5557
$ws.Range("K132").Value = 8983.667099999999
$ws.Range("M132").Value = -6453.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 47725
$ws.Range("J98").Value = 47725
$ws.Range("L98").Value = 47725
$ws.Range("N98").Value = -53715
$ws.Range("H132").Value = 1838.2593
$ws.Range("I132").Value = 1384.25
$ws.Range("K132").Value = 4152.75
$ws.Range("M132").Value = -1622.75

